# Update the FHIR ValueSet metadata on the "Metadata" worksheet and the
# "Include from American Dental " worksheet to reflect the move from the
# Alvearie/IBM project to LinuxForHealth.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/ada-tooth"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

$wsCodes = $wb.Worksheets.Item("Include from American Dental ")
$wsCodes.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/ada-tooth"
